# Updated with Field Name with standard name and Minor release
#
# Changes applied:
#  1. Both sheets ("API-Testing" and "API-Testing-Sheet2-Duplicate"):
#     - Delete column F (RequestProcessingType) entirely - it shifts every
#       column from G onward one to the left.
#     - Rename header cells to their "standard" names:
#         HTTPAction      -> Action
#         ExcludeField    -> ExcludeFields
#         HttpStatusCode  -> StatusCode
#         security        -> Security
#         tags            -> Tags
#  2. Sheet1 column J (StatusCode) gets an explicit width.
#  3. Selection / active cell bookkeeping + active tab (minor release, second
#     sheet was left active) to match the "Minor release" part of the commit.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # --- structural edit: drop the RequestProcessingType column -----------
    $ws.Columns.Item(6).Delete()

    # --- rename headers to their standardized names -----------------------
    $ws.Cells.Item(1, 8).Value  = "Action"
    $ws.Cells.Item(1, 9).Value  = "ExcludeFields"
    $ws.Cells.Item(1, 10).Value = "StatusCode"
    $ws.Cells.Item(1, 13).Value = "Security"
    $ws.Cells.Item(1, 14).Value = "Tags"
}

# Sheet-specific column width / selection bookkeeping
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# StatusCode column on sheet1 was resized
$ws1.Columns.Item(10).ColumnWidth = 13.6

# restore per-sheet selections
$ws1.Range("J8").Select()
$ws2.Range("O1").Select()

# Minor release: the workbook was left with the second sheet active
$ws2.Activate()
